# Rubrica Andamento - apply diff changes via Excel COM interop
$wb = $excel.ActiveWorkbook

# 1. Rename existing sheet "Sheet1" -> "Escopo"
$wb.Worksheets.Item("Sheet1").Name = "Escopo"

# 2. Add new blank sheet named "RN", then move it after "Escopo"
$wsRN = $wb.Worksheets.Add()
$wsRN.Name = "RN"
$wsRN.Move($null, $wb.Worksheets.Item("Escopo"))

# Select cell E13 on RN sheet (matches end-state selection in diff)
$wb.Worksheets.Item("RN").Range("E13").Select()

# 3. Update values on Escopo sheet
$ws = $wb.Worksheets.Item("Escopo")

# Row 2: E2 changes from "Em Andamento" to "Concluido"; F2 0.75 -> 1
$ws.Range("E2").Value = "Concluido"
$ws.Range("F2").Value = 1

# Row 4: E4 gets "Em Andamento" (was empty); F4 0 -> 0.2
$ws.Range("E4").Value = "Em Andamento"
$ws.Range("F4").Value = 0.2

# Row 5: F5 0.15 -> 0.3 (E5 remains "Em Andamento")
$ws.Range("F5").Value = 0.3

# Set active selection to E7 on Escopo sheet (matches end-state diff)
$ws.Activate()
$ws.Range("E7").Select()
